# ---------------------------------------------------------------------------
# This presentation ships two theme parts:
#   ppt/theme/theme1.xml  - bound to the (only) Slide Master / the design
#                            actually used by every slide. Currently the
#                            "Integral" / "Red Violet" palette.
#   ppt/theme/theme2.xml  - bound only to the Notes Master. Currently the
#                            stock "Office Theme" palette.
#
# The target edit swaps the two palettes: the Slide Master's theme becomes
# the "Office Theme" colours (what theme2.xml used to hold) while the Notes
# Master's theme becomes the old "Integral" colours.  Font scheme and format
# scheme are identical between the two themes already, so only the 12
# colour-scheme slots (background/text/accents/hyperlinks) actually change.
#
# We drive this the same way a user would from the Design tab / the Slide
# Master's colour-scheme editor: via the ColorScheme object exposed on the
# presentation's master, updating each of the 12 theme colour slots in turn
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - PowerPoint's classic
# 1-based ColorScheme.Colors(index) ordering) to the new "Office Theme" RGB
# values.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# New palette for the presentation's (Slide Master) theme -- this is the
# palette theme2.xml/"Office Theme" currently uses.
$officeTheme = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501    # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

$masterScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Colors($i).RGB = $officeTheme[$i]
}

Write-Host "Theme colour scheme updated (Slide Master -> Office Theme)."
